$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay a text value even though it looks numeric
    # (mirrors the source data, which stores every figure as inline text).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# --- Simple price updates (column D) ---
Set-TextValue "D2"  "243.74"
Set-TextValue "D3"  "23.22"
Set-TextValue "D4"  "5.416"
Set-TextValue "D5"  "0.05971"
Set-TextValue "D6"  "3.428"
Set-TextValue "D7"  "6.530"
Set-TextValue "D8"  "0.8121"
Set-TextValue "D9"  "0.9359"
Set-TextValue "D10" "0.1423"
Set-TextValue "D11" "0.07425"
Set-TextValue "D12" "0.03285"
Set-TextValue "D13" "0.03068"
Set-TextValue "D14" "0.09350"
Set-TextValue "D15" "3.854"
Set-TextValue "D16" "0.001574"
Set-TextValue "D17" "0.04704"
Set-TextValue "D18" "0.0005901"
Set-TextValue "D19" "0.005897"
Set-TextValue "D20" "0.001275"
Set-TextValue "D21" "0.004911"
Set-TextValue "D22" "0.00006800"
Set-TextValue "D27" "0.0002340"
Set-TextValue "D40" "0.03970"

# --- Rows 41-43: BKEXToken, CEJI and KickToken get re-ranked/rotated ---
# Row 41 becomes CEJI
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D41" "0.005201"
$ws.Range("E41").Value = "40CEJICEJIBestin24h"

# Row 42 becomes KickToken
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.006440"
$ws.Range("E42").Value = "41KickTokenKICK"

# Row 43 becomes BKEXToken
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1080"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Row 44: LocalTraders price + label change ---
Set-TextValue "D44" "0.009197"
$ws.Range("E44").Value = "43LocalTradersLCT"

# --- Row 45: CoinLion price ---
Set-TextValue "D45" "0.00005215"

# --- Row 47: CoinbaseStockToken label gains "Worstin24h" ---
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# --- Row 48: BOLO price ---
Set-TextValue "D48" "0.002308"
